# Apply cryptos list update (scraped values refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so Excel keeps exact
# text (e.g. trailing zeros) instead of coercing the value to a number
$textCells = @("D5", "D6", "D7", "D8", "D9", "D12", "D13", "D14", "D17", "D20", "D21", "D22", "D23", "D25", "D27", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D41", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update coin metadata + price/volume figures
$ws.Range("D2").Value = "73.931.86"
$ws.Range("E2").Value = "  +8.59%  "
$ws.Range("D3").Value = "2.585.01"
$ws.Range("E3").Value = "  +7.07%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "582.36"
$ws.Range("E5").Value = "  +4.93%  "
$ws.Range("D6").Value = "181.04"
$ws.Range("E6").Value = "  +13.95%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  +5.04%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.205"
$ws.Range("E9").Value = "  +26.18%  "
$ws.Range("D10").Value = "2.587.40"
$ws.Range("E10").Value = "  +7.28%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "0.354"
$ws.Range("E12").Value = "  +7.23%  "
$ws.Range("D13").Value = "4.78"
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("D14").Value = "0.0000193"
$ws.Range("E14").Value = "  +11.33%  "
$ws.Range("D15").Value = "73.656.98"
$ws.Range("E15").Value = "  +8.40%  "
$ws.Range("D16").Value = "3.050.78"
$ws.Range("E16").Value = "  +6.94%  "
$ws.Range("D17").Value = "26.01"
$ws.Range("E17").Value = "  +13.97%  "
$ws.Range("D18").Value = "2.565.17"
$ws.Range("E18").Value = "  +6.42%  "
$ws.Range("E19").Value = "  +11.94%  "
$ws.Range("D20").Value = "7.91"
$ws.Range("E20").Value = "  +15.32%  "
$ws.Range("D21").Value = "364.26"
$ws.Range("E21").Value = "  +9.95%  "
$ws.Range("D22").Value = "2.23"
$ws.Range("E22").Value = "  +18.26%  "
$ws.Range("D23").Value = "4.05"
$ws.Range("E23").Value = "  +6.71%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "69.00"
$ws.Range("E25").Value = "  +4.26%  "
$ws.Range("E26").Value = "  +12.85%  "
$ws.Range("D27").Value = "9.07"
$ws.Range("E27").Value = "  +10.99%  "
$ws.Range("D28").Value = "2.709.31"
$ws.Range("E28").Value = "  +6.86%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("D30").Value = "0.0₃0935"
$ws.Range("E30").Value = "  +15.71%  "
$ws.Range("D31").Value = "7.87"
$ws.Range("E31").Value = "  +10.97%  "
$ws.Range("D32").Value = "496.41"
$ws.Range("E32").Value = "  +17.69%  "
$ws.Range("D33").Value = "1.33"
$ws.Range("E33").Value = "  +17.34%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").Value = "1.69"
$ws.Range("E35").Value = "  +5.88%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.118"
$ws.Range("E36").Value = "  +13.01%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "160.17"
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("D38").Value = "19.06"
$ws.Range("E38").Value = "  +6.95%  "
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "4.82"
$ws.Range("E41").Value = "  +12.31%  "
$ws.Range("E42").Value = "  +12.53%  "
$ws.Range("D43").Value = "0.317"
$ws.Range("E43").Value = "  +7.24%  "
$ws.Range("E44").Value = "  +20.72%  "
$ws.Range("D45").Value = "39.26"
$ws.Range("E45").Value = "  +5.43%  "
$ws.Range("D46").Value = "1.14"
$ws.Range("E46").Value = "  +6.71%  "
$ws.Range("D47").Value = "147.50"
$ws.Range("E47").Value = "  +11.47%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.0795"
$ws.Range("E48").Value = "  +11.65%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "3.56"
$ws.Range("E49").Value = "  +7.34%  "
$ws.Range("D50").Value = "0.517"
$ws.Range("E50").Value = "  +7.97%  "
$ws.Range("D51").Value = "0.582"
$ws.Range("E51").Value = "  +5.05%  "
